# Update basic ignored file for ISMIP6 and RFMIP (step 3) due to drq version 01.00.27.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Force plain text (avoids "1" -> numeric auto-conversion) without leaving
    # a residual custom number-format style on the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).NumberFormat = "General"
}

function Set-HyperlinkCell($addr, $url) {
    $ws.Range($addr).Formula = '=HYPERLINK("' + $url + '","web")'
}

function Copy-IceSheetCommentStyle($addr) {
    # Rows 440 / 442 already carry the "We do not have an Antarctic ice
    # sheet." comment in the Cambria-font style (cellXf index 2) - reuse it.
    $ws.Range("H440").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Row 454 - IyrAnt / modelCellAreai (ISMIP6)
# ---------------------------------------------------------------------
$ws.Range("A454").Value = "IyrAnt"
$ws.Range("B454").Value = "modelCellAreai"
Set-TextCell "C454" "1"
$ws.Range("D454").Value = "xant yant time"
$ws.Range("E454").Value = "The cell area of the ice sheet model."
$ws.Range("F454").Value = "m2"
Set-HyperlinkCell "G454" "http://clipc-services.ceda.ac.uk/dreq/u/865d0e00-53e6-11e6-b524-5404a60d96b5.html"
Copy-IceSheetCommentStyle "H454"
$ws.Range("H454").Value = "We do not have an Antarctic ice sheet."
$ws.Range("I454").Value = "Thomas"
$ws.Range("J454").Value = "Horizontal area of ice-sheet grid cells"
$ws.Range("K454").Value = "ISMIP6"

# ---------------------------------------------------------------------
# Row 455 - IyrAnt / sftgif (ISMIP6)
# ---------------------------------------------------------------------
$ws.Range("A455").Value = "IyrAnt"
$ws.Range("B455").Value = "sftgif"
Set-TextCell "C455" "1"
$ws.Range("D455").Value = "xant yant time typeli"
$ws.Range("E455").Value = "Fraction of Grid Cell Covered with Glacier"
$ws.Range("F455").Value = "%"
Set-HyperlinkCell "G455" "http://clipc-services.ceda.ac.uk/dreq/u/a1d2e309c6f25017442ad6c79c4f9eca.html"
Copy-IceSheetCommentStyle "H455"
$ws.Range("H455").Value = "We do not have an Antarctic ice sheet."
$ws.Range("I455").Value = "Thomas"
$ws.Range("J455").Value = "Fraction of grid cell covered by land ice (ice sheet, ice shelf, ice cap, glacier)"
$ws.Range("K455").Value = "ISMIP6"

# ---------------------------------------------------------------------
# Row 456 - IyrAnt / sftgrf (ISMIP6)
# ---------------------------------------------------------------------
$ws.Range("A456").Value = "IyrAnt"
$ws.Range("B456").Value = "sftgrf"
Set-TextCell "C456" "1"
$ws.Range("D456").Value = "xant yant time typegis"
$ws.Range("E456").Value = "Grounded Ice Sheet  Area Fraction"
$ws.Range("F456").Value = "%"
Set-HyperlinkCell "G456" "http://clipc-services.ceda.ac.uk/dreq/u/590e5de4-9e49-11e5-803c-0d0b866b59f3.html"
Copy-IceSheetCommentStyle "H456"
$ws.Range("H456").Value = "We do not have an Antarctic ice sheet."
$ws.Range("I456").Value = "Thomas"
$ws.Range("J456").Value = "Fraction of grid cell covered by grounded ice sheet"
$ws.Range("K456").Value = "ISMIP6"

# Row 457 intentionally stays blank (gap in the source sheet).

# ---------------------------------------------------------------------
# Row 458 - Efx / rlu (RFMIP)
# ---------------------------------------------------------------------
$ws.Range("A458").Value = "Efx"
$ws.Range("B458").Value = "rlu"
Set-TextCell "C458" "1"
$ws.Range("D458").Value = "alevhalf spectband"
$ws.Range("E458").Value = "Upwelling Longwave Radiation"
$ws.Range("F458").Value = "W m-2"
Set-HyperlinkCell "G458" "http://clipc-services.ceda.ac.uk/dreq/u/bcfeacf77d49ef51a6ee66a1ab0ebcb4.html"
$ws.Range("H458").Value = "Not available in IFS: All Up and downwelling radiation is only at the TOA and the surface available in IFS standard output. In IFS it is not possible to distinguish output in spectral intervals. Note here also global area and time averages are asked. Or maybe output at a certain diagnostic time step is meant? (No grib code available on table 128 -  Grib 1 for different spectral bands). Would it be possible to output 2 spectral bands: UV and NIR both of them diffuse and parallel. We need to ask expert of the radiation code."
$ws.Range("I458").Value = "Twan & Thomas"
$ws.Range("J458").Value = "Upwelling longwave radiation (includes the fluxes at the surface and TOA)"
$ws.Range("K458").Value = "RFMIP"

# ---------------------------------------------------------------------
# Row 459 - Efx / rsu (RFMIP)
# ---------------------------------------------------------------------
$ws.Range("A459").Value = "Efx"
$ws.Range("B459").Value = "rsu"
Set-TextCell "C459" "1"
$ws.Range("D459").Value = "alevhalf spectband"
$ws.Range("E459").Value = "Upwelling Shortwave Radiation"
$ws.Range("F459").Value = "W m-2"
Set-HyperlinkCell "G459" "http://clipc-services.ceda.ac.uk/dreq/u/c323f38340e4846931ad4891232d839d.html"
$ws.Range("H459").Value = "Not available in IFS: All Up and downwelling radiation is only at the TOA and the surface available in IFS standard output. In IFS it is not possible to distinguish output in spectral intervals. Note here also global area and time averages are asked. Or maybe output at a certain diagnostic time step is meant? (No grib code available on table 128 -  Grib 1 for different spectral bands). Would it be possible to output 2 spectral bands: UV and NIR both of them diffuse and parallel. We need to ask expert of the radiation code."
$ws.Range("I459").Value = "Twan & Thomas"
$ws.Range("J459").Value = "Upwelling shortwave radiation  (includes also the fluxes at the surface and top of atmosphere)"
$ws.Range("K459").Value = "RFMIP"

# ---------------------------------------------------------------------
# Row 460 - Efx / rld (RFMIP)
# ---------------------------------------------------------------------
$ws.Range("A460").Value = "Efx"
$ws.Range("B460").Value = "rld"
Set-TextCell "C460" "1"
$ws.Range("D460").Value = "alevhalf spectband"
$ws.Range("E460").Value = "Downwelling Longwave Radiation"
$ws.Range("F460").Value = "W m-2"
Set-HyperlinkCell "G460" "http://clipc-services.ceda.ac.uk/dreq/u/c432bfbfc0e7f4403f91af39736ff61c.html"
$ws.Range("H460").Value = "Not available in IFS: All Up and downwelling radiation is only at the TOA and the surface available in IFS standard output. In IFS it is not possible to distinguish output in spectral intervals. Note here also global area and time averages are asked. Or maybe output at a certain diagnostic time step is meant? (No grib code available on table 128 -  Grib 1 for different spectral bands). Would it be possible to output 2 spectral bands: UV and NIR both of them diffuse and parallel. We need to ask expert of the radiation code."
$ws.Range("I460").Value = "Twan & Thomas"
$ws.Range("J460").Value = "Downwelling Longwave Radiation (includes the fluxes at the surface and TOA)"
$ws.Range("K460").Value = "RFMIP"

# ---------------------------------------------------------------------
# Row 461 - Efx / rsd (RFMIP)
# ---------------------------------------------------------------------
$ws.Range("A461").Value = "Efx"
$ws.Range("B461").Value = "rsd"
Set-TextCell "C461" "1"
$ws.Range("D461").Value = "alevhalf spectband"
$ws.Range("E461").Value = "Downwelling Shortwave Radiation"
$ws.Range("F461").Value = "W m-2"
Set-HyperlinkCell "G461" "http://clipc-services.ceda.ac.uk/dreq/u/eb9ac643cd9c73cae960d6d2db7b901d.html"
$ws.Range("H461").Value = "Not available in IFS: All Up and downwelling radiation is only at the TOA and the surface available in IFS standard output. In IFS it is not possible to distinguish output in spectral intervals. Note here also global area and time averages are asked. Or maybe output at a certain diagnostic time step is meant? (No grib code available on table 128 -  Grib 1 for different spectral bands). Would it be possible to output 2 spectral bands: UV and NIR both of them diffuse and parallel. We need to ask expert of the radiation code."
$ws.Range("I461").Value = "Twan & Thomas"
$ws.Range("J461").Value = "Downwelling shortwave radiation (includes the fluxes at the surface and top-of-atmosphere)"
$ws.Range("K461").Value = "RFMIP"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Cosmetic side-effects of the original edit session: the sheet's used
# range grew to cover the full sheet height (as happens after the
# author scrolled/selected near the very bottom of the sheet) and a
# handful of trailing blank rows picked up explicit row heights.
# Reproduce the resulting navigation/view state.
# ---------------------------------------------------------------------
$ws.Range("A1048572:A1048576").RowHeight = 12.8
$ws.Cells.Item(1048576, 11).Font.Size = 11
$ws.Range("A458").Select()
